$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# "Repeats the code <expression> times. If a variable ..." paragraph:
# insert the word "local " so it reads "If a local variable ...".
# ---------------------------------------------------------------------

$full = $d.Content.Text
$idx = $full.IndexOf("If a variable")
$insPoint = $idx + ("If a ").Length
$insRange = $d.Range($insPoint, $insPoint)
$insRange.InsertBefore("local ")

# Inserting text merges every run that shares identical formatting with the
# edit point into a single run (this also swallows the neighbouring
# unrelated runs that just happen to carry the same rPr, e.g. the smart
# quotes around "index"). Re-split that merged block back into the runs we
# actually want by toggling a cosmetic property across each boundary - this
# forces a run split without altering the visible formatting.
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf("If a local ")
$p0 = $idx2 + ("If a ").Length          # start of "local "
$p1 = $p0 + ("local ").Length           # start of "variable " (end of "local ")
$p2 = $p1 + ("variable ").Length        # start of the opening smart quote
$p3 = $p2 + 1                           # start of "index"
$p4 = $p3 + ("index").Length            # start of the closing smart quote
$p5 = $p4 + 1                           # start of " exists ... at the "

function Split-Run([int]$s, [int]$e) {
    if ($e -gt $s) {
        $r = $d.Range($s, $e)
        $r.Bold = 1
        $r.Bold = 0
    }
}

Split-Run $p0 $p1
Split-Run $p1 $p2
Split-Run $p2 $p3
Split-Run $p3 $p4
Split-Run $p4 $p5

# Word keeps the "_GoBack" bookmark at the site of the most recent edit.
# Re-adding it here (right after "local ", before "variable") moves it from
# wherever it previously sat (after the "Passing variable addresses ..."
# paragraph) to this new insertion point.
$bkRange = $d.Range($p1, $p1)
$d.Bookmarks.Add("_GoBack", $bkRange)
